$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking row (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update Total row (B12): 27 -> 45
$ws.Range("B12").Value = 45

# Update the correct/total marks label (E12): "19/84" -> "45/140"
$ws.Range("E12").Value = "45/140"
